$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 576 (shifts the existing 576:629 block down to 579:632)
$ws.Rows("576:578").Insert()

# New weekly data block (row 576 = "Especial", 577 = "Primera", 578 = "Segunda")
$newRows = @(
    @{ Row = 576; Quality = "Especial"; Date = 45106; Volumen = 300; Min = 24000; Max = 25000; Prom = 24500; PrecioKg = 3500 },
    @{ Row = 577; Quality = "Primera";  Date = 45106; Volumen = 400; Min = 20000; Max = 21000; Prom = 20500; PrecioKg = 2929 },
    @{ Row = 578; Quality = "Segunda";  Date = 45106; Volumen = 360; Min = 15000; Max = 16000; Prom = 15500; PrecioKg = 2214 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 2
    $ws.Cells.Item($row, 2).Value = "Comercializadora del Agro de Limarí"
    $ws.Cells.Item($row, 3).Value = "Coquimbo"
    $ws.Cells.Item($row, 4).Value = $r.Date
    $ws.Cells.Item($row, 5).Value = 4
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100101
    $ws.Cells.Item($row, 8).Value = "Berries"
    $ws.Cells.Item($row, 9).Value = 100112025
    $ws.Cells.Item($row, 10).Value = "Frutilla"
    $ws.Cells.Item($row, 11).Value = "Sin especificar"
    $ws.Cells.Item($row, 12).Value = $r.Quality
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.Min
    $ws.Cells.Item($row, 15).Value = $r.Max
    $ws.Cells.Item($row, 16).Value = $r.Prom
    $ws.Cells.Item($row, 17).Value = "$/bandeja 7 kilos"
    $ws.Cells.Item($row, 18).Value = "Provincia de Melipilla"
    $ws.Cells.Item($row, 19).Value = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value = 7
}
